$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark that sits right after the
#    "SkyNet" run (it is a hidden bookmark, addressable by name even
#    though it does not show up in Bookmarks.Count).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Add a new paragraph right after the first empty paragraph that
#    follows "SkyNet", carrying the programmer credit line.
$skyNetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "SkyNet") {
        $skyNetPara = $d.Paragraphs.Item($i)
        break
    }
}

$blankAfterSkyNet = $skyNetPara.Next()
$blankAfterSkyNet.Range.InsertParagraphAfter()

$newPara = $blankAfterSkyNet.Next()
$iAcute = [char]0x00ED
$newPara.Range.Text = "Programador 1: Oscar Joel Delc" + $iAcute + "d Revolorio 2017173: Modelo ER y Modelo Entidad"
